# TestWorkbook.xlsx (ReadTableFromExcel regression data) update:
# Add tests for ReadTableFromExcel for numeric column headings and to
# better handle include and exclude filters.
#
# Adds two new worksheets:
#   - "NumericHeadings" (copy of Address_NoFormulas, header row replaced
#     with plain numbers instead of Column#_Type strings)
#   - "ExcludeMultiple"  (copy of Address_NoFormulas, two of the String
#     values in column D tweaked so include/exclude filters have near
#     duplicates to exercise)
# and positions them in the sheet tab order, updating each sheet's
# selection / active tab to match.
#
# NOTE: worksheet references are re-fetched via $wb.Worksheets.Item(name)
# right before each use (rather than reusing an earlier PowerShell
# variable) because a sheet reference can go stale across an operation
# that mutates the sheet collection (Copy / Move / rename).

$wb = $excel.ActiveWorkbook

$srcSheetName = "Address_NoFormulas"

# --- Create "NumericHeadings" first (copy appended at the end) ------------
$src = $wb.Worksheets.Item($srcSheetName)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$src.Copy($null, $lastSheet)
$wb.Worksheets.Item($wb.Worksheets.Count).Name = "NumericHeadings"

# --- Then create "ExcludeMultiple" (copy appended at the end) -------------
$src = $wb.Worksheets.Item($srcSheetName)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$src.Copy($null, $lastSheet)
$wb.Worksheets.Item($wb.Worksheets.Count).Name = "ExcludeMultiple"

# --- Move the new sheets into their final tab positions --------------------
# Final order: Address_NoFormulas, ExcludeMultiple, Address_Formulas,
#              NumericHeadings, Worksheet_NoHeaders, Worksheet_Headers
$wb.Worksheets.Item("ExcludeMultiple").Move($wb.Worksheets.Item("Address_Formulas"))
$wb.Worksheets.Item("NumericHeadings").Move($wb.Worksheets.Item("Worksheet_NoHeaders"))

# --- Populate "ExcludeMultiple" data tweaks --------------------------------
$wb.Worksheets.Item("ExcludeMultiple").Range("D7").Value = "aString 6"
$wb.Worksheets.Item("ExcludeMultiple").Range("D11").Value = "bString 10"

# --- Populate "NumericHeadings" header row with numbers --------------------
$wsNumeric = $wb.Worksheets.Item("NumericHeadings")
$wsNumeric.Range("B1").Value = 1990
$wsNumeric.Range("C1").Value = 1991
$wsNumeric.Range("D1").Value = 1992
$wsNumeric.Range("E1").Value = 1993.3
$wsNumeric.Range("F1").Value = 1994
$wsNumeric.Range("G1").Value = 1995
$wsNumeric.Range("E1").NumberFormat = "0"
$wsNumeric.Range("F1").NumberFormat = "0"

# --- Selections per sheet (match the edited workbook) ----------------------
[void]$wb.Worksheets.Item("Address_NoFormulas").Rows(9).Select()
[void]$wb.Worksheets.Item("Address_Formulas").Range("E12").Select()
[void]$wb.Worksheets.Item("Worksheet_NoHeaders").Range("C15").Select()
[void]$wb.Worksheets.Item("Worksheet_Headers").Range("D26").Select()
[void]$wb.Worksheets.Item("NumericHeadings").Range("G2").Select()
[void]$wb.Worksheets.Item("ExcludeMultiple").Range("D14").Select()

# "ExcludeMultiple" is the active/visible tab.
$wb.Worksheets.Item("ExcludeMultiple").Activate()
